$d = $word.ActiveDocument

$replacements = @(
    @{old="69×62="; new="61×56="},
    @{old="85×16="; new="77×59="},
    @{old="34×96="; new="99×42="},
    @{old="71×87="; new="20×95="},
    @{old="72×79="; new="90×63="},
    @{old="86×56="; new="29×82="},
    @{old="64×16="; new="21×93="},
    @{old="48×19="; new="17×87="},
    @{old="52×73="; new="44×98="},
    @{old="70×45="; new="71×56="},
    @{old="11×92="; new="12×55="},
    @{old="23×93="; new="33×58="},
    @{old="93×27="; new="31×91="},
    @{old="70×53="; new="27×92="},
    @{old="97×95="; new="12×58="},
    @{old="97×89="; new="82×31="},
    @{old="37×95="; new="32×91="},
    @{old="97×73="; new="49×69="},
    @{old="16×31="; new="75×70="},
    @{old="82×77="; new="30×19="},
    @{old="33×41="; new="35×76="},
    @{old="29×15="; new="99×90="},
    @{old="45×15="; new="46×71="},
    @{old="62×47="; new="98×68="},
    @{old="50×38="; new="89×85="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
